$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking values are not
# auto-converted into actual numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.314.29"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.483.37"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "520.27"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").Value = "134.90"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.559"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "2.500.40"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "0.0989"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "0.340"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "2.926.04"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "58.264.67"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "22.15"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "2.491.46"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "321.13"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("D23").Value = "5.76"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "64.46"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "0.0₃0751"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "169.85"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "6.32"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "18.14"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "36.67"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("D41").Value = "0.799"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "5.21"
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("D43").Value = "278.14"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "3.46"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "124.10"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("D47").Value = "0.0910"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "0.0492"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "0.0214"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "17.10"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "1.741.10"
$ws.Range("E51").Value = "  -0.60%  "
